$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = "A"
$ws.Range("E3").Value = "G"
$ws.Range("F3").Value = -0.3653988847128944

$ws.Range("D8").Value = "A"
$ws.Range("E8").Value = "G"
$ws.Range("F8").Value = 0.2593599748104452

$ws.Range("D9").Value = "C"
$ws.Range("E9").Value = "T"
$ws.Range("F9").Value = -0.3513041745623368
